# cryptos.xlsx periodic data refresh (GitHub Actions bot).
#
# Updates the Price (col D) and Volume(1h) (col E) columns with the
# latest scraped quotes, and corrects the Polygon/Polkadot row order
# (rows 15-16 had been swapped).
#
# Column D sometimes holds values that look like numbers (e.g. "228.43",
# "0.380"). The sheet stores these as literal text (so formats like a
# trailing zero or thousand-dot-separated big numbers like "38.261.86"
# survive untouched), so a leading apostrophe is used to force Excel to
# keep the entry as text instead of auto-converting it to a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.261.86'
$ws.Range('E2').Value = '  +1.54%  '

$ws.Range('D3').Value = '2.092.72'
$ws.Range('E3').Value = '  +3.22%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '''228.43'
$ws.Range('E5').Value = '  +0.49%  '

$ws.Range('E6').Value = '  +0.88%  '

$ws.Range('D7').Value = '''60.92'
$ws.Range('E7').Value = '  +1.53%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('D9').Value = '''0.380'
$ws.Range('E9').Value = '  +1.38%  '

$ws.Range('D10').Value = '''0.0843'
$ws.Range('E10').Value = '  +1.66%  '

$ws.Range('E11').Value = '  +0.23%  '

$ws.Range('D12').Value = '2.405.08'
$ws.Range('E12').Value = '  +3.27%  '

$ws.Range('D13').Value = '''14.74'
$ws.Range('E13').Value = '  +2.23%  '

$ws.Range('D14').Value = '''22.27'
$ws.Range('E14').Value = '  +5.89%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '''5.45'
$ws.Range('E15').Value = '  +5.40%  '

$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '''0.775'
$ws.Range('E16').Value = '  +2.25%  '

$ws.Range('D17').Value = '2.092.35'
$ws.Range('E17').Value = '  +3.72%  '

$ws.Range('D18').Value = '38.244.51'
$ws.Range('E18').Value = '  +1.52%  '

$ws.Range('D19').Value = '''70.37'
$ws.Range('E19').Value = '  +1.38%  '

$ws.Range('E20').Value = '  +1.52%  '

$ws.Range('D21').Value = '0.0₃0832'
$ws.Range('E21').Value = '  +1.12%  '

$ws.Range('D22').Value = '''224.93'
$ws.Range('E22').Value = '  +0.75%  '

$ws.Range('E23').Value = '  +0.01%  '

$ws.Range('D24').Value = '''2.43'
$ws.Range('E24').Value = '  +1.92%  '

$ws.Range('D25').Value = '''2.31'
$ws.Range('E25').Value = '  +3.10%  '

$ws.Range('D26').Value = '''170.06'
$ws.Range('E26').Value = '  +1.55%  '

$ws.Range('D27').Value = '''9.39'
$ws.Range('E27').Value = '  +1.65%  '

$ws.Range('E28').Value = '  +1.40%  '

$ws.Range('D29').Value = '''18.98'
$ws.Range('E29').Value = '  +1.08%  '

$ws.Range('E30').Value = '  +8.87%  '

$ws.Range('E31').Value = '  -0.30%  '

$ws.Range('D32').Value = '''2.34'
$ws.Range('E32').Value = '  +5.28%  '

$ws.Range('D33').Value = '''4.73'
$ws.Range('E33').Value = '  +6.07%  '

$ws.Range('D34').Value = '''4.43'
$ws.Range('E34').Value = '  +1.14%  '

$ws.Range('D35').Value = '''0.0603'
$ws.Range('E35').Value = '  +0.11%  '

$ws.Range('E36').Value = '  +1.35%  '

$ws.Range('D37').Value = '''2.38'
$ws.Range('E37').Value = '  +4.30%  '

$ws.Range('D38').Value = '''3.50'
$ws.Range('E38').Value = '  +5.18%  '

$ws.Range('E39').Value = '  +0.15%  '

$ws.Range('D40').Value = '''18.19'
$ws.Range('E40').Value = '  +1.96%  '

$ws.Range('D41').Value = '1.543.10'
$ws.Range('E41').Value = '  +0.38%  '

$ws.Range('D42').Value = '''99.76'
$ws.Range('E42').Value = '  +4.57%  '

$ws.Range('D43').Value = '''0.0219'
$ws.Range('E43').Value = '  +1.71%  '

$ws.Range('E44').Value = '  +1.36%  '

$ws.Range('D45').Value = '''0.0908'
$ws.Range('E45').Value = '  -0.22%  '

$ws.Range('E46').Value = '  +0.24%  '

$ws.Range('D47').Value = '''1.11'
$ws.Range('E47').Value = '  +1.25%  '

$ws.Range('D48').Value = '''7.50'
$ws.Range('E48').Value = '  +5.62%  '

$ws.Range('E49').Value = '  +2.67%  '

$ws.Range('E50').Value = '  +0.71%  '

$ws.Range('D51').Value = '2.292.03'
$ws.Range('E51').Value = '  +3.29%  '
